# Translate the two column headers from German to English
# (shared-string table: "Ventilstellung" -> "travel", "Heizleistung" -> "heating power")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "travel"
$ws.Range("B1").Value = "heating power"

# Move / leave the sheet's active selection on B1 (matches saved <selection> in sheetView)
$ws.Range("B1").Select()
